$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Type Pattern" sheet: update the two example cells (C12, C13) that
#    describe "I fill blanket with a value" examples. The plain-comma
#    phrasing is replaced by a sentence-ending punctuation + "E.g." on a
#    new line, keeping the trailing number (red) as its own run.
# ---------------------------------------------------------------------
$wsType = $wb.Worksheets.Item("Type Pattern")

# C12: full-width period variant, value stays "10km"
$wsType.Range("C12").Value = "I fill blanket with a value。`nE.g. 10km"
$wsType.Range("C12").Characters(35, 4).Font.Color = 255
$wsType.Range("C12").WrapText = $true

# C13: regular period variant, value stays "10 km"
$wsType.Range("C13").Value = "I fill blanket with a value.`nE.g. 10 km"
$wsType.Range("C13").Characters(35, 5).Font.Color = 255
$wsType.Range("C13").WrapText = $true

# ---------------------------------------------------------------------
# 2. Active sheet / selection moves from "other pattern" to
#    "Type Pattern", with a new selected range and scroll position.
# ---------------------------------------------------------------------
$wsType.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsType.Range("B10:D13").Select()
